# Automatische test-sync: 2025-06-26 22:25:50
# Adds a new logged e-mail (row 25) to the "Logs" sheet and updates the
# "Dashboard" sheet category counts/order to reflect it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- Append the new log entry on row 25 -----------------------------------
$ws.Cells.Item(25, 1).Value = "Wat zijn de verzendkosten?"
$ws.Cells.Item(25, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item(25, 3).Value = "Testmail #2: Wat zijn de verzendkosten?"
$ws.Cells.Item(25, 4).Value = "Productinformatie"
$ws.Cells.Item(25, 5).Value = "Beste afzender,`nDank u voor uw interesse in onze producten/diensten. Om u nauwkeurige informatie te verstrekken over de verzendkosten, hebben we meer details nodig zoals het product/dienst waar u naar informeert en het afleveradres. Zou u ons kunnen voorzien van deze informatie zodat we u een precieze schatting van de verzendkosten kunnen geven?`nMet vriendelijke groet,`n[Naam]`nE-mailassistent"
$ws.Cells.Item(25, 6).Value = "2025-06-26 22:24:51"
$ws.Cells.Item(25, 7).Value = "Ja"
$ws.Cells.Item(25, 8).Value = "Nee"
$ws.Cells.Item(25, 9).Value = "Ja"

# The multi-line answer in column E triggers an automatic row-height /
# word-wrap recalculation; re-running AutoFit keeps the row on the
# worksheet's default height, matching the rest of the sheet.
$ws.Rows.Item(25).AutoFit()

# --- Extend the conditional formatting ranges to cover the new row --------
function Extend-ConditionalFormatting($column) {
    $conditions = $ws.Range($column + "2").FormatConditions
    for ($i = 1; $i -le $conditions.Count; $i++) {
        $conditions.Item($i).ModifyAppliesToRange($ws.Range($column + "2:" + $column + "25"))
    }
}
Extend-ConditionalFormatting "D"
Extend-ConditionalFormatting "G"
Extend-ConditionalFormatting "H"
Extend-ConditionalFormatting "I"

# --- Update the Dashboard sheet with the recalculated category counts -----
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("A3").Value = "Productinformatie"
$dashboard.Range("A5").Value = "Openingstijden / Locatie"
$dashboard.Range("B5").Value = 2
$dashboard.Range("A6").Value = "Retour / Terugbetaling"
